$d = $word.ActiveDocument

# Remove the "_GoBack" bookmark that splits the bold run
# "DOCX, DOC, PDF, HTML, XPS, R" / "TF and TXT" into two runs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Merge the now-adjacent bold runs "DOCX, DOC, PDF, HTML, XPS, R" and "TF and TXT"
# into a single run with text "DOCX, DOC, PDF, HTML, XPS, RTF and TXT".
# Step 1: replace with a distinct placeholder so Word consolidates the match into one run.
$d.Content.Find.Execute("DOCX, DOC, PDF, HTML, XPS, RTF and TXT", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GEMBOX_PLACEHOLDER_TEXT", 2)

# Step 2: replace the placeholder with the final desired text (now a single run).
$d.Content.Find.Execute("GEMBOX_PLACEHOLDER_TEXT", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2)
